# 自动更新Excel文件
# For every data row, advance one day: "剩余" (days remaining) decreases by 1.
# When a row's remaining count has hit 1 (i.e. would expire), renew it instead:
# reset "剩余" back to the full "总天" count and roll "开始时间" forward by
# the renewal window (10 days) from its previous start date.
# Rows whose "开始时间" isn't a clean 8-digit yyyyMMdd value can't be parsed
# into a date for the renewal calc, so (matching the historical run) they are
# left untouched and simply skipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    try {
        $total = $ws.Range("D$r").Value2
        $remain = $ws.Range("E$r").Value2
        $start = $ws.Range("F$r").Value2

        if ($null -eq $remain) { continue }

        $startStr = [string]([int]$start)
        if ($startStr.Length -ne 8) {
            throw "invalid start date '$startStr' on row $r"
        }

        if ($remain -le 1) {
            # About to run out -> renew the cycle.
            $startDate = [DateTime]::ParseExact($startStr, "yyyyMMdd", $null)
            $newStartDate = $startDate.AddDays(10)
            $newStart = [int]$newStartDate.ToString("yyyyMMdd")

            $ws.Range("E$r").Value = $total
            $ws.Range("F$r").Value = $newStart
        }
        else {
            $ws.Range("E$r").Value = $remain - 1
        }
    }
    catch {
        # Unparseable row (e.g. a corrupted date) - leave it as-is, same as
        # the original automation would skip it on error.
        continue
    }
}
